$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column (D) cells are treated as text, matching the
# original inline-string content (many values look numeric, e.g. "0.9994")
$ws.Range("D2:D51").NumberFormat = "@"

# Row 2
$ws.Range("D2").Value = "29.209.85"
$ws.Range("E2").Value = "  -0.97%  "

# Row 3
$ws.Range("D3").Value = "1.859.55"
$ws.Range("E3").Value = "  -0.61%  "

# Row 4
$ws.Range("D4").Value = "0.9994"
$ws.Range("E4").Value = "  -0.19%  "

# Row 5
$ws.Range("D5").Value = "0.7155"
$ws.Range("E5").Value = "  -0.33%  "

# Row 6
$ws.Range("D6").Value = "240.51"
$ws.Range("E6").Value = "  +0.54%  "

# Row 7
$ws.Range("D7").Value = "0.9999"
$ws.Range("E7").Value = "  -0.12%  "

# Row 8
$ws.Range("D8").Value = "0.07753"
$ws.Range("E8").Value = "  -0.80%  "

# Row 9
$ws.Range("D9").Value = "0.3078"
$ws.Range("E9").Value = "  +0.26%  "

# Row 10
$ws.Range("D10").Value = "25.20"
$ws.Range("E10").Value = "  -0.36%  "

# Row 11
$ws.Range("D11").Value = "0.08255"
$ws.Range("E11").Value = "  +0.24%  "

# Row 12
$ws.Range("D12").Value = "5.237"
$ws.Range("E12").Value = "  +0.16%  "

# Row 13
$ws.Range("B13").Value = "Polygon"
$ws.Range("C13").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D13").Value = "0.7176"
$ws.Range("E13").Value = "  -0.52%  "

# Row 14
$ws.Range("B14").Value = "WrappedEther"
$ws.Range("C14").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D14").Value = "1.842.31"
$ws.Range("E14").Value = "  -2.10%  "

# Row 15
$ws.Range("D15").Value = "90.26"
$ws.Range("E15").Value = "  +0.30%  "

# Row 16
$ws.Range("D16").Value = "29.195.37"
$ws.Range("E16").Value = "  -1.19%  "

# Row 17
$ws.Range("D17").Value = "5.870"
$ws.Range("E17").Value = "  +0.74%  "

# Row 18
$ws.Range("D18").Value = "244.15"
$ws.Range("E18").Value = "  +1.34%  "

# Row 19
$ws.Range("D19").Value = "0.000007799"
$ws.Range("E19").Value = "  -0.62%  "

# Row 20
$ws.Range("D20").Value = "13.15"
$ws.Range("E20").Value = "  -1.12%  "

# Row 21
$ws.Range("D21").Value = "2.110.51"
$ws.Range("E21").Value = "  -0.76%  "

# Row 22
$ws.Range("E22").Value = "  -0.07%  "

# Row 23
$ws.Range("D23").Value = "7.972"
$ws.Range("E23").Value = "  +3.10%  "

# Row 24
$ws.Range("D24").Value = "0.9992"
$ws.Range("E24").Value = "  -0.25%  "

# Row 25
$ws.Range("D25").Value = "0.1596"
$ws.Range("E25").Value = "  +1.86%  "

# Row 26
$ws.Range("D26").Value = "162.46"
$ws.Range("E26").Value = "  -0.17%  "

# Row 27
$ws.Range("D27").Value = "8.940"
$ws.Range("E27").Value = "  -0.23%  "

# Row 28
$ws.Range("D28").Value = "18.28"
$ws.Range("E28").Value = "  -0.21%  "

# Row 29
$ws.Range("E29").Value = "  +0.94%  "

# Row 30
$ws.Range("E30").Value = "  -3.22%  "

# Row 31
$ws.Range("E31").Value = "  +1.76%  "

# Row 32
$ws.Range("D32").Value = "4.195"
$ws.Range("E32").Value = "  +2.93%  "

# Row 33
$ws.Range("E33").Value = "  -1.06%  "

# Row 34
$ws.Range("D34").Value = "1.909"
$ws.Range("E34").Value = "  -1.31%  "

# Row 35
$ws.Range("E35").Value = "  -2.10%  "

# Row 36
$ws.Range("D36").Value = "0.7270"
$ws.Range("E36").Value = "  +1.52%  "

# Row 37
$ws.Range("D37").Value = "2.675"
$ws.Range("E37").Value = "  +0.08%  "

# Row 38
$ws.Range("E38").Value = "  -0.55%  "

# Row 39
$ws.Range("D39").Value = "2.685"
$ws.Range("E39").Value = "  -1.32%  "

# Row 40
$ws.Range("D40").Value = "1.148.75"
$ws.Range("E40").Value = "  -2.11%  "

# Row 41
$ws.Range("D41").Value = "0.9060"
$ws.Range("E41").Value = "  -0.13%  "

# Row 42
$ws.Range("D42").Value = "6.146"
$ws.Range("E42").Value = "  +2.48%  "

# Row 43
$ws.Range("D43").Value = "72.26"
$ws.Range("E43").Value = "  +1.10%  "

# Row 44
$ws.Range("D44").Value = "0.9997"
$ws.Range("E44").Value = "  -0.16%  "

# Row 45
$ws.Range("D45").Value = "101.66"
$ws.Range("E45").Value = "  -0.61%  "

# Row 46
$ws.Range("D46").Value = "2.003.35"
$ws.Range("E46").Value = "  -1.75%  "

# Row 47
$ws.Range("D47").Value = "0.5227"
$ws.Range("E47").Value = "  -2.52%  "

# Row 48
$ws.Range("D48").Value = "1.769"
$ws.Range("E48").Value = "  +0.16%  "

# Row 49
$ws.Range("B49").Value = "EnergySwap"
$ws.Range("C49").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D49").Value = "9.332"
$ws.Range("E49").Value = "  +1.90%  "

# Row 50
$ws.Range("B50").Value = "BabyDogeCoin"
$ws.Range("C50").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D50").Value = "0.00000000119"
$ws.Range("E50").Value = "  -1.41%  "

# Row 51
$ws.Range("E51").Value = "  +1.46%  "
